# Addition of Yes/No Button to UX Testing Flow
# Update the "Wait time (sec) after Video ends" value (cell E2) from 90s to 22s,
# and move the active selection to E3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E2 holds a time-of-day formatted number representing seconds-as-a-fraction-of-a-day.
# 22 seconds = 22/86400 of a day.
$ws.Range("E2").Value = 22/86400

# Move the selection down to E3, matching the new active cell in the workbook.
$ws.Range("E3").Select()
